$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 799
$ws.Range("I18").Value = 799
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 799
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -515

$ws.Range("H40").Value = 3790.25
$ws.Range("I40").Value = 3307.4285
$ws.Range("J40").Value = 4165.778
$ws.Range("K40").Value = 3307.4285
$ws.Range("L40").Value = 4165.778
$ws.Range("M40").Value = -3132.4285
$ws.Range("N40").Value = -4515.778

$ws.Range("H74").Value = 8683.8125
$ws.Range("I74").Value = 6793.1113
$ws.Range("K74").Value = 6793.1113
$ws.Range("M74").Value = -5857.1113

$ws.Range("H77").Value = 8683.8125
$ws.Range("I77").Value = 6793.1113
$ws.Range("K77").Value = 33965.5565
$ws.Range("M77").Value = -29285.5565

$ws.Range("H98").Value = 6314.3335
$ws.Range("I98").Value = 505.81818
$ws.Range("K98").Value = 505.81818
$ws.Range("M98").Value = 992.18182

$ws.Range("H122").Value = 6314.3335
$ws.Range("I122").Value = 505.81818
$ws.Range("K122").Value = 1517.45454
$ws.Range("M122").Value = 932.54546

$ws.Range("H132").Value = 1152.7858
$ws.Range("I132").Value = 1132.1951
$ws.Range("K132").Value = 3396.5853
$ws.Range("M132").Value = -866.5852999999997

$ws.Range("H135").Value = 613.6667
$ws.Range("I135").Value = 652.875
$ws.Range("J135").Value = 488.2
$ws.Range("K135").Value = 5875.875
$ws.Range("L135").Value = 4393.8
$ws.Range("M135").Value = -3340.875
$ws.Range("N135").Value = -9463.799999999999

$ws.Range("H137").Value = 3002.3684
$ws.Range("I137").Value = 2435.125
$ws.Range("K137").Value = 7305.375
$ws.Range("M137").Value = -4755.375

$ws.Range("H138").Value = 4968.5386
$ws.Range("I138").Value = 4819
$ws.Range("J138").Value = 5062
$ws.Range("K138").Value = 14457
$ws.Range("L138").Value = 15186
$ws.Range("M138").Value = -9317
$ws.Range("N138").Value = -25466

$ws.Range("H141").Value = 573.36365
$ws.Range("I141").Value = 573.36365
$ws.Range("K141").Value = 1720.09095
$ws.Range("M141").Value = 3459.90905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3904
$ws.Range("I61").Value = 4116.353
$ws.Range("K61").Value = 4116.353
$ws.Range("M61").Value = -3904.353

$ws.Range("H63").Value = 6999.875
$ws.Range("I63").Value = 2999.5
$ws.Range("K63").Value = 2999.5
$ws.Range("M63").Value = -2313.5

$ws.Range("H66").Value = 6999.875
$ws.Range("I66").Value = 2999.5
$ws.Range("K66").Value = 14997.5
$ws.Range("M66").Value = -11565.5

$ws.Range("H74").Value = 10418028
$ws.Range("I74").Value = 11495497
$ws.Range("K74").Value = 11495497
$ws.Range("M74").Value = -11494623

$ws.Range("H77").Value = 10418028
$ws.Range("I77").Value = 11495497
$ws.Range("K77").Value = 57477485
$ws.Range("M77").Value = -57473117

$ws.Range("H110").Value = 4054.9285
$ws.Range("I110").Value = 1474.3
$ws.Range("J110").Value = 10506.5
$ws.Range("K110").Value = 1474.3
$ws.Range("L110").Value = 10506.5
$ws.Range("M110").Value = 570.7
$ws.Range("N110").Value = -14596.5

$ws.Range("H132").Value = 4119.8066
$ws.Range("I132").Value = 2862.52
$ws.Range("J132").Value = 9358.5
$ws.Range("K132").Value = 8587.559999999999
$ws.Range("L132").Value = 28075.5
$ws.Range("M132").Value = -6057.559999999999
$ws.Range("N132").Value = -33135.5

$ws.Range("H136").Value = 3904
$ws.Range("I136").Value = 4116.353
$ws.Range("K136").Value = 12349.059
$ws.Range("M136").Value = -9799.059000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2383.625
$ws.Range("I80").Value = 1998
$ws.Range("J80").Value = 2615
$ws.Range("K80").Value = 1998
$ws.Range("L80").Value = 2615
$ws.Range("M80").Value = -1000
$ws.Range("N80").Value = -4611

$ws.Range("H83").Value = 2383.625
$ws.Range("I83").Value = 1998
$ws.Range("J83").Value = 2615
$ws.Range("K83").Value = 9990
$ws.Range("L83").Value = 13075
$ws.Range("M83").Value = -4998
$ws.Range("N83").Value = -23059

$ws.Range("H99").Value = 2950
$ws.Range("I99").Value = 2950
$ws.Range("K99").Value = 2950
$ws.Range("M99").Value = -1452

$ws.Range("H134").Value = 2024.9678
$ws.Range("I134").Value = 1133.7858
$ws.Range("K134").Value = 3401.3574
$ws.Range("M134").Value = -866.3574000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29287.756
$ws.Range("I31").Value = 3634.5625
$ws.Range("J31").Value = 120499.11
$ws.Range("K31").Value = 3634.5625
$ws.Range("L31").Value = 120499.11
$ws.Range("M31").Value = -3339.5625
$ws.Range("N31").Value = -121089.11

$ws.Range("H34").Value = 29287.756
$ws.Range("I34").Value = 3634.5625
$ws.Range("J34").Value = 120499.11
$ws.Range("K34").Value = 3634.5625
$ws.Range("L34").Value = 120499.11
$ws.Range("M34").Value = -3432.5625
$ws.Range("N34").Value = -120903.11

$ws.Range("H59").Value = 89127
$ws.Range("J59").Value = 89127
$ws.Range("L59").Value = 89127
$ws.Range("N59").Value = -91417

$ws.Range("H103").Value = 35749.5
$ws.Range("I103").Value = 35749.5
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 35749.5
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -34577.5

$ws.Range("H132").Value = 1972.6571
$ws.Range("I132").Value = 1230
$ws.Range("J132").Value = 5562.1665
$ws.Range("K132").Value = 3690
$ws.Range("L132").Value = 16686.4995
$ws.Range("M132").Value = -1160
$ws.Range("N132").Value = -21746.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 302
$ws.Range("I11").Value = 1000
$ws.Range("K11").Value = 3000
$ws.Range("M11").Value = -2860

$ws.Range("H114").Value = 462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7695.727
$ws.Range("I80").Value = 6449.857
$ws.Range("J80").Value = 9876
$ws.Range("K80").Value = 6449.857
$ws.Range("L80").Value = 9876
$ws.Range("M80").Value = -5451.857
$ws.Range("N80").Value = -11872

$ws.Range("H83").Value = 7695.727
$ws.Range("I83").Value = 6449.857
$ws.Range("J83").Value = 9876
$ws.Range("K83").Value = 32249.285
$ws.Range("L83").Value = 49380
$ws.Range("M83").Value = -27257.285
$ws.Range("N83").Value = -59364

$ws.Range("H122").Value = 15005.538
$ws.Range("I122").Value = 15606.6
$ws.Range("K122").Value = 46819.8
$ws.Range("M122").Value = -44369.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6805.375
$ws.Range("I40").Value = 5872.5186
$ws.Range("K40").Value = 5872.5186
$ws.Range("M40").Value = -5736.5186

$ws.Range("H132").Value = 3678.6775
$ws.Range("I132").Value = 2324.182
$ws.Range("J132").Value = 6989.6665
$ws.Range("K132").Value = 6972.545999999999
$ws.Range("L132").Value = 20968.9995
$ws.Range("M132").Value = -4442.545999999999
$ws.Range("N132").Value = -26028.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1320.8889
$ws.Range("I107").Value = 1007.8
$ws.Range("J107").Value = 1712.25
$ws.Range("K107").Value = 3023.4
$ws.Range("L107").Value = 5136.75
$ws.Range("M107").Value = -1103.4
$ws.Range("N107").Value = -8976.75

$ws.Range("H122").Value = 15977.55
$ws.Range("J122").Value = 17659.75
$ws.Range("L122").Value = 52979.25
$ws.Range("N122").Value = -57879.25

$ws.Range("H132").Value = 6917.0557
$ws.Range("I132").Value = 4348.4814
$ws.Range("K132").Value = 13045.4442
$ws.Range("M132").Value = -10515.4442

$ws.Range("H136").Value = 1975
$ws.Range("I136").Value = 1474.0385
$ws.Range("K136").Value = 4422.1155
$ws.Range("M136").Value = -1872.1155
